$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain Text (they look numeric, e.g. "231.34")
# by pre-setting NumberFormat to Text before assigning the value.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.713.36'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.803.96'
$ws.Range('E3').Value = '  -1.35%  '
$ws.Range('E4').Value = '  +0.59%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.34'
$ws.Range('E5').Value = '  -2.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5941'
$ws.Range('E6').Value = '  -2.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.005'
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2781'
$ws.Range('E8').Value = '  -1.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06850'
$ws.Range('E9').Value = '  -3.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.42'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07549'
$ws.Range('E11').Value = '  -1.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.799.46'
$ws.Range('E12').Value = '  -1.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.729'
$ws.Range('E13').Value = '  -2.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6275'
$ws.Range('E14').Value = '  -1.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.048.62'
$ws.Range('E15').Value = '  -1.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009290'
$ws.Range('E16').Value = '  -7.98%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '75.46'
$ws.Range('E17').Value = '  -4.88%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '28.692.05'
$ws.Range('E18').Value = '  -1.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.482'
$ws.Range('E19').Value = '  -7.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.005'
$ws.Range('E20').Value = '  +0.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '211.55'
$ws.Range('E21').Value = '  -7.34%  '
$ws.Range('E22').Value = '  -3.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.868'
$ws.Range('E23').Value = '  -2.33%  '
$ws.Range('E24').Value = '  +0.53%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.29'
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.842'
$ws.Range('E26').Value = '  -2.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1271'
$ws.Range('E27').Value = '  -2.00%  '
$ws.Range('E28').Value = '  -1.06%  '
$ws.Range('E29').Value = '  -2.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06178'
$ws.Range('E30').Value = '  -5.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.420'
$ws.Range('E31').Value = '  -2.66%  '
$ws.Range('E32').Value = '  -1.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.755'
$ws.Range('E33').Value = '  -1.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.716'
$ws.Range('E34').Value = '  -1.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.056'
$ws.Range('E35').Value = '  -6.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6419'
$ws.Range('E36').Value = '  -1.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.505'
$ws.Range('E37').Value = '  -1.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.725'
$ws.Range('E38').Value = '  -1.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01712'
$ws.Range('E39').Value = '  -2.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.428'
$ws.Range('E40').Value = '  -1.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.142.03'
$ws.Range('E41').Value = '  -6.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8662'
$ws.Range('E42').Value = '  -7.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.005'
$ws.Range('E43').Value = '  +0.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.81'
$ws.Range('E44').Value = '  -0.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.963.54'
$ws.Range('E45').Value = '  -0.96%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '60.56'
$ws.Range('E46').Value = '  -4.23%  '
$ws.Range('E47').Value = '  -5.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.588'
$ws.Range('E48').Value = '  -1.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.376'
$ws.Range('E49').Value = '  -2.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05472'
$ws.Range('E50').Value = '  -1.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4502'
$ws.Range('E51').Value = '  -1.31%  '
